# Update the Home Loan amortization table with the latest projected figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserInfo")

# Make sure the percentage-looking text values are stored as plain text,
# not auto-converted to numeric percentages by Excel.
$pctCells = "G6","G7","G8","G9","G10","G11","G12","G13","G14","G15"
foreach ($addr in $pctCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 6 - 2024
$ws.Range("B6").Value = "₹ 1,71,591"
$ws.Range("C6").Value = "₹ 2,35,038"
$ws.Range("D6").Value = "₹ 35,000"
$ws.Range("E6").Value = "₹ 4,41,629"
$ws.Range("F6").Value = "₹ 30,38,409"
$ws.Range("G6").Value = "5.35%"

# Row 7 - 2025
$ws.Range("B7").Value = "₹ 2,23,571"
$ws.Range("C7").Value = "₹ 2,64,384"
$ws.Range("F7").Value = "₹ 28,14,838"
$ws.Range("G7").Value = "12.31%"

# Row 8 - 2026
$ws.Range("B8").Value = "₹ 2,44,544"
$ws.Range("C8").Value = "₹ 2,43,411"
$ws.Range("F8").Value = "₹ 25,70,294"
$ws.Range("G8").Value = "19.93%"

# Row 9 - 2027
$ws.Range("B9").Value = "₹ 2,67,484"
$ws.Range("C9").Value = "₹ 2,20,471"
$ws.Range("F9").Value = "₹ 23,02,810"
$ws.Range("G9").Value = "28.26%"

# Row 10 - 2028
$ws.Range("B10").Value = "₹ 2,92,576"
$ws.Range("C10").Value = "₹ 1,95,380"
$ws.Range("F10").Value = "₹ 20,10,234"
$ws.Range("G10").Value = "37.38%"

# Row 11 - 2029
$ws.Range("B11").Value = "₹ 3,20,021"
$ws.Range("C11").Value = "₹ 1,67,934"
$ws.Range("F11").Value = "₹ 16,90,213"
$ws.Range("G11").Value = "47.35%"

# Row 12 - 2030
$ws.Range("B12").Value = "₹ 3,50,041"
$ws.Range("C12").Value = "₹ 1,37,914"
$ws.Range("F12").Value = "₹ 13,40,172"
$ws.Range("G12").Value = "58.25%"

# Row 13 - 2031
$ws.Range("B13").Value = "₹ 3,82,878"
$ws.Range("C13").Value = "₹ 1,05,077"
$ws.Range("F13").Value = "₹ 9,57,294"
$ws.Range("G13").Value = "70.18%"

# Row 14 - 2032
$ws.Range("B14").Value = "₹ 4,18,794"
$ws.Range("C14").Value = "₹ 69,161"
$ws.Range("F14").Value = "₹ 5,38,500"
$ws.Range("G14").Value = "83.22%"

# Row 15 - 2033
$ws.Range("B15").Value = "₹ 4,58,080"
$ws.Range("C15").Value = "₹ 29,875"
$ws.Range("F15").Value = "₹ 80,420"
$ws.Range("G15").Value = "97.49%"

# Row 16 - 2034 (B16 shares its value with F15's "Balance" figure)
$ws.Range("B16").Value = "₹ 80,420"
$ws.Range("C16").Value = "₹ 906"
$ws.Range("D16").Value = "₹ 7,000"
$ws.Range("E16").Value = "₹ 88,326"
